$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectList")

# Add the officer NRIC value for the "Acacia Breeze" project (row 2, column N - "Officer")
$ws.Range("N2").Value = "T1234567J"
